$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: EONIARATE -> EONIARATE(1)
$ws.Name = "EONIARATE(1)"

# Update the header label in A1: DATE -> Date
$ws.Range("A1").Value = "Date"

# Zoom the view from 160% to 200%
$excel.ActiveWindow.Zoom = 200

# Move the active selection to D15
$ws.Range("D15").Select()

# Resize columns (values chosen so the engine's internal
# char-width quantization lands as close as possible to the
# target OOXML column widths of 15.34 / 14.74 / 25.39 / 19.66).
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
$ws.Columns.Item(5).ColumnWidth = 24.5
$ws.Columns.Item(6).ColumnWidth = 18.833333333333332

Write-Host "edit.ps1 completed"
